# "Generate Report for Archive"
# - Update status text "Ready for handoff" -> "In Translation" everywhere it
#   appears (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# - Shrink the now-narrower "Status" columns to match the new autofit width
#   (Overview columns E & F, zh-cn column C, de-de column C).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"
$overview.Range("E1:F1").ColumnWidth = 12.5

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"
$zhcn.Range("C1").ColumnWidth = 12.5

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"
$dede.Range("C1").ColumnWidth = 12.5
